$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Weekly data refresh ---
# Existing price records (rows 2-67) are updated with refreshed values
# (shifted dates / prices / origins), and a new record is appended as row 68.

# Row 2
$ws.Cells.Item(2, 4).Value = 44230
$ws.Cells.Item(2, 11).Value = 22000
$ws.Cells.Item(2, 12).Value = 24000
$ws.Cells.Item(2, 13).Value = 23000
$ws.Cells.Item(2, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(2, 15).Value = 'Región del Maule'
$ws.Cells.Item(2, 16).Value = 920

# Row 3
$ws.Cells.Item(3, 4).Value = 44937
$ws.Cells.Item(3, 8).Value = 'Sin especificar'
$ws.Cells.Item(3, 11).Value = 38000
$ws.Cells.Item(3, 12).Value = 40000
$ws.Cells.Item(3, 13).Value = 39000
$ws.Cells.Item(3, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(3, 16).Value = 1560

# Row 4
$ws.Cells.Item(4, 4).Value = 44568
$ws.Cells.Item(4, 10).Value = 200
$ws.Cells.Item(4, 11).Value = 25000
$ws.Cells.Item(4, 12).Value = 26000
$ws.Cells.Item(4, 13).Value = 25500
$ws.Cells.Item(4, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(4, 16).Value = 1020

# Row 5
$ws.Cells.Item(5, 4).Value = 44651
$ws.Cells.Item(5, 10).Value = 140
$ws.Cells.Item(5, 11).Value = 20000
$ws.Cells.Item(5, 12).Value = 23000
$ws.Cells.Item(5, 13).Value = 21714
$ws.Cells.Item(5, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(5, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(5, 16).Value = 869

# Row 6
$ws.Cells.Item(6, 4).Value = 44624
$ws.Cells.Item(6, 10).Value = 110
$ws.Cells.Item(6, 11).Value = 26000
$ws.Cells.Item(6, 12).Value = 27000
$ws.Cells.Item(6, 13).Value = 26545
$ws.Cells.Item(6, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(6, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(6, 16).Value = 1062

# Row 7
$ws.Cells.Item(7, 4).Value = 44574
$ws.Cells.Item(7, 8).Value = 'Magnum'
$ws.Cells.Item(7, 10).Value = 100
$ws.Cells.Item(7, 11).Value = 27000
$ws.Cells.Item(7, 12).Value = 28000
$ws.Cells.Item(7, 13).Value = 27500
$ws.Cells.Item(7, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(7, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(7, 16).Value = 1100

# Row 8
$ws.Cells.Item(8, 4).Value = 44638
$ws.Cells.Item(8, 10).Value = 180
$ws.Cells.Item(8, 11).Value = 23000
$ws.Cells.Item(8, 12).Value = 24000
$ws.Cells.Item(8, 13).Value = 23444
$ws.Cells.Item(8, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(8, 16).Value = 938

# Row 9
$ws.Cells.Item(9, 4).Value = 44321
$ws.Cells.Item(9, 11).Value = 24000
$ws.Cells.Item(9, 12).Value = 25000
$ws.Cells.Item(9, 13).Value = 24500
$ws.Cells.Item(9, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(9, 15).Value = 'Región del Maule'
$ws.Cells.Item(9, 16).Value = 980

# Row 10
$ws.Cells.Item(10, 4).Value = 44572
$ws.Cells.Item(10, 10).Value = 250
$ws.Cells.Item(10, 11).Value = 26000
$ws.Cells.Item(10, 12).Value = 27000
$ws.Cells.Item(10, 13).Value = 26480
$ws.Cells.Item(10, 16).Value = 1059

# Row 11
$ws.Cells.Item(11, 4).Value = 44265
$ws.Cells.Item(11, 11).Value = 20000
$ws.Cells.Item(11, 12).Value = 22000
$ws.Cells.Item(11, 13).Value = 21000
$ws.Cells.Item(11, 16).Value = 840

# Row 12
$ws.Cells.Item(12, 4).Value = 44441
$ws.Cells.Item(12, 11).Value = 28000
$ws.Cells.Item(12, 12).Value = 29000
$ws.Cells.Item(12, 13).Value = 28500
$ws.Cells.Item(12, 16).Value = 1140

# Row 13
$ws.Cells.Item(13, 4).Value = 44769
$ws.Cells.Item(13, 11).Value = 34000
$ws.Cells.Item(13, 12).Value = 35000
$ws.Cells.Item(13, 13).Value = 34500
$ws.Cells.Item(13, 16).Value = 1380

# Row 14
$ws.Cells.Item(14, 4).Value = 44783
$ws.Cells.Item(14, 11).Value = 38000
$ws.Cells.Item(14, 12).Value = 40000
$ws.Cells.Item(14, 13).Value = 39000
$ws.Cells.Item(14, 16).Value = 1560

# Row 15
$ws.Cells.Item(15, 4).Value = 44811
$ws.Cells.Item(15, 11).Value = 27000
$ws.Cells.Item(15, 12).Value = 28000
$ws.Cells.Item(15, 13).Value = 27500
$ws.Cells.Item(15, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(15, 15).Value = 'Perú'
$ws.Cells.Item(15, 16).Value = 1100

# Row 16
$ws.Cells.Item(16, 4).Value = 44532
$ws.Cells.Item(16, 8).Value = 'Magnum'
$ws.Cells.Item(16, 10).Value = 250
$ws.Cells.Item(16, 11).Value = 33000
$ws.Cells.Item(16, 12).Value = 35000
$ws.Cells.Item(16, 13).Value = 33800
$ws.Cells.Item(16, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(16, 16).Value = 1352

# Row 17
$ws.Cells.Item(17, 4).Value = 44279
$ws.Cells.Item(17, 11).Value = 28000
$ws.Cells.Item(17, 12).Value = 30000
$ws.Cells.Item(17, 13).Value = 29000
$ws.Cells.Item(17, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(17, 15).Value = 'Región del Maule'
$ws.Cells.Item(17, 16).Value = 1160

# Row 18
$ws.Cells.Item(18, 4).Value = 44468
$ws.Cells.Item(18, 8).Value = 'Sin especificar'
$ws.Cells.Item(18, 11).Value = 31000
$ws.Cells.Item(18, 12).Value = 32000
$ws.Cells.Item(18, 13).Value = 31500
$ws.Cells.Item(18, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(18, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(18, 16).Value = 1260

# Row 19
$ws.Cells.Item(19, 4).Value = 44945
$ws.Cells.Item(19, 10).Value = 250
$ws.Cells.Item(19, 11).Value = 14000
$ws.Cells.Item(19, 12).Value = 15000
$ws.Cells.Item(19, 13).Value = 14600
$ws.Cells.Item(19, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(19, 15).Value = 'Región del Maule'
$ws.Cells.Item(19, 16).Value = 584

# Row 20
$ws.Cells.Item(20, 4).Value = 44797
$ws.Cells.Item(20, 10).Value = 100
$ws.Cells.Item(20, 11).Value = 44000
$ws.Cells.Item(20, 12).Value = 45000
$ws.Cells.Item(20, 13).Value = 44500
$ws.Cells.Item(20, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(20, 15).Value = 'Perú'
$ws.Cells.Item(20, 16).Value = 1780

# Row 21
$ws.Cells.Item(21, 4).Value = 44664
$ws.Cells.Item(21, 11).Value = 22000
$ws.Cells.Item(21, 12).Value = 24000
$ws.Cells.Item(21, 13).Value = 23000
$ws.Cells.Item(21, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(21, 16).Value = 920

# Row 22
$ws.Cells.Item(22, 4).Value = 44363
$ws.Cells.Item(22, 10).Value = 100
$ws.Cells.Item(22, 11).Value = 25000
$ws.Cells.Item(22, 12).Value = 26000
$ws.Cells.Item(22, 13).Value = 25500
$ws.Cells.Item(22, 15).Value = 'Perú'
$ws.Cells.Item(22, 16).Value = 1020

# Row 23
$ws.Cells.Item(23, 4).Value = 44461
$ws.Cells.Item(23, 11).Value = 33000
$ws.Cells.Item(23, 12).Value = 34000
$ws.Cells.Item(23, 13).Value = 33500
$ws.Cells.Item(23, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(23, 16).Value = 1340

# Row 24
$ws.Cells.Item(24, 4).Value = 44706
$ws.Cells.Item(24, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(24, 15).Value = 'Perú'

# Row 25
$ws.Cells.Item(25, 4).Value = 44253
$ws.Cells.Item(25, 8).Value = 'Magnum'
$ws.Cells.Item(25, 10).Value = 200
$ws.Cells.Item(25, 11).Value = 25000
$ws.Cells.Item(25, 12).Value = 26000
$ws.Cells.Item(25, 13).Value = 25500
$ws.Cells.Item(25, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(25, 16).Value = 1020

# Row 26
$ws.Cells.Item(26, 4).Value = 44876
$ws.Cells.Item(26, 10).Value = 70
$ws.Cells.Item(26, 11).Value = 32000
$ws.Cells.Item(26, 12).Value = 34000
$ws.Cells.Item(26, 13).Value = 32857
$ws.Cells.Item(26, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(26, 15).Value = 'Perú'
$ws.Cells.Item(26, 16).Value = 1314

# Row 27
$ws.Cells.Item(27, 4).Value = 44678
$ws.Cells.Item(27, 11).Value = 19000
$ws.Cells.Item(27, 12).Value = 20000
$ws.Cells.Item(27, 13).Value = 19500
$ws.Cells.Item(27, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(27, 16).Value = 780

# Row 28
$ws.Cells.Item(28, 4).Value = 44237
$ws.Cells.Item(28, 8).Value = 'Sin especificar'
$ws.Cells.Item(28, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(28, 15).Value = 'Región del Maule'

# Row 29
$ws.Cells.Item(29, 4).Value = 44603
$ws.Cells.Item(29, 10).Value = 120
$ws.Cells.Item(29, 11).Value = 25000
$ws.Cells.Item(29, 12).Value = 26000
$ws.Cells.Item(29, 13).Value = 25583
$ws.Cells.Item(29, 16).Value = 1023

# Row 30
$ws.Cells.Item(30, 4).Value = 44323
$ws.Cells.Item(30, 8).Value = 'Magnum'
$ws.Cells.Item(30, 11).Value = 20000
$ws.Cells.Item(30, 12).Value = 22000
$ws.Cells.Item(30, 13).Value = 21000
$ws.Cells.Item(30, 15).Value = 'Perú'
$ws.Cells.Item(30, 16).Value = 840

# Row 31
$ws.Cells.Item(31, 4).Value = 44622
$ws.Cells.Item(31, 10).Value = 220
$ws.Cells.Item(31, 11).Value = 24000
$ws.Cells.Item(31, 12).Value = 25000
$ws.Cells.Item(31, 13).Value = 24545
$ws.Cells.Item(31, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(31, 16).Value = 982

# Row 32
$ws.Cells.Item(32, 4).Value = 44609
$ws.Cells.Item(32, 10).Value = 200
$ws.Cells.Item(32, 11).Value = 28000
$ws.Cells.Item(32, 12).Value = 30000
$ws.Cells.Item(32, 13).Value = 29000
$ws.Cells.Item(32, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(32, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(32, 16).Value = 1160

# Row 33
$ws.Cells.Item(33, 4).Value = 44335
$ws.Cells.Item(33, 8).Value = 'Magnum'
$ws.Cells.Item(33, 10).Value = 100
$ws.Cells.Item(33, 11).Value = 35000
$ws.Cells.Item(33, 12).Value = 36000
$ws.Cells.Item(33, 13).Value = 35500
$ws.Cells.Item(33, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(33, 16).Value = 1420

# Row 34
$ws.Cells.Item(34, 4).Value = 44435
$ws.Cells.Item(34, 10).Value = 100
$ws.Cells.Item(34, 11).Value = 25000
$ws.Cells.Item(34, 12).Value = 26000
$ws.Cells.Item(34, 13).Value = 25500
$ws.Cells.Item(34, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(34, 15).Value = 'Perú'
$ws.Cells.Item(34, 16).Value = 1020

# Row 35
$ws.Cells.Item(35, 4).Value = 44615
$ws.Cells.Item(35, 8).Value = 'Sin especificar'
$ws.Cells.Item(35, 10).Value = 100
$ws.Cells.Item(35, 11).Value = 28000
$ws.Cells.Item(35, 12).Value = 30000
$ws.Cells.Item(35, 13).Value = 29000
$ws.Cells.Item(35, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(35, 15).Value = 'Región del Maule'
$ws.Cells.Item(35, 16).Value = 1160

# Row 36
$ws.Cells.Item(36, 4).Value = 44594
$ws.Cells.Item(36, 10).Value = 200
$ws.Cells.Item(36, 11).Value = 34000
$ws.Cells.Item(36, 12).Value = 35000
$ws.Cells.Item(36, 13).Value = 34500
$ws.Cells.Item(36, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(36, 15).Value = 'Región del Maule'
$ws.Cells.Item(36, 16).Value = 1380

# Row 37
$ws.Cells.Item(37, 4).Value = 44530
$ws.Cells.Item(37, 8).Value = 'Sin especificar'
$ws.Cells.Item(37, 10).Value = 110
$ws.Cells.Item(37, 11).Value = 22000
$ws.Cells.Item(37, 12).Value = 23000
$ws.Cells.Item(37, 13).Value = 22455
$ws.Cells.Item(37, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(37, 15).Value = 'Región del Maule'
$ws.Cells.Item(37, 16).Value = 898

# Row 38
$ws.Cells.Item(38, 4).Value = 44188
$ws.Cells.Item(38, 10).Value = 100
$ws.Cells.Item(38, 11).Value = 38000
$ws.Cells.Item(38, 12).Value = 40000
$ws.Cells.Item(38, 13).Value = 39000
$ws.Cells.Item(38, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(38, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(38, 16).Value = 1560

# Row 39
$ws.Cells.Item(39, 4).Value = 44708
$ws.Cells.Item(39, 11).Value = 20000
$ws.Cells.Item(39, 12).Value = 22000
$ws.Cells.Item(39, 13).Value = 21000
$ws.Cells.Item(39, 16).Value = 840

# Row 40
$ws.Cells.Item(40, 4).Value = 44203
$ws.Cells.Item(40, 10).Value = 100
$ws.Cells.Item(40, 11).Value = 20000
$ws.Cells.Item(40, 12).Value = 22000
$ws.Cells.Item(40, 13).Value = 21000
$ws.Cells.Item(40, 15).Value = 'Región del Maule'
$ws.Cells.Item(40, 16).Value = 840

# Row 41
$ws.Cells.Item(41, 4).Value = 44167
$ws.Cells.Item(41, 8).Value = 'Sin especificar'
$ws.Cells.Item(41, 11).Value = 18000
$ws.Cells.Item(41, 12).Value = 19000
$ws.Cells.Item(41, 13).Value = 18500
$ws.Cells.Item(41, 15).Value = 'Región del Maule'
$ws.Cells.Item(41, 16).Value = 740

# Row 42
$ws.Cells.Item(42, 4).Value = 44868
$ws.Cells.Item(42, 8).Value = 'Magnum'
$ws.Cells.Item(42, 11).Value = 27000
$ws.Cells.Item(42, 12).Value = 28000
$ws.Cells.Item(42, 13).Value = 27500
$ws.Cells.Item(42, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(42, 15).Value = 'Perú'
$ws.Cells.Item(42, 16).Value = 1100

# Row 43
$ws.Cells.Item(43, 4).Value = 44580
$ws.Cells.Item(43, 8).Value = 'Magnum'
$ws.Cells.Item(43, 15).Value = 'Región Metropolitana'

# Row 44
$ws.Cells.Item(44, 4).Value = 44672
$ws.Cells.Item(44, 8).Value = 'Sin especificar'
$ws.Cells.Item(44, 10).Value = 130
$ws.Cells.Item(44, 11).Value = 24000
$ws.Cells.Item(44, 12).Value = 25000
$ws.Cells.Item(44, 13).Value = 24615
$ws.Cells.Item(44, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(44, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(44, 16).Value = 985

# Row 45
$ws.Cells.Item(45, 4).Value = 44475
$ws.Cells.Item(45, 11).Value = 44000
$ws.Cells.Item(45, 12).Value = 45000
$ws.Cells.Item(45, 13).Value = 44500
$ws.Cells.Item(45, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(45, 15).Value = 'Perú'
$ws.Cells.Item(45, 16).Value = 1780

# Row 46
$ws.Cells.Item(46, 4).Value = 44447
$ws.Cells.Item(46, 8).Value = 'Magnum'
$ws.Cells.Item(46, 11).Value = 37000
$ws.Cells.Item(46, 12).Value = 38000
$ws.Cells.Item(46, 13).Value = 37500
$ws.Cells.Item(46, 15).Value = 'Perú'
$ws.Cells.Item(46, 16).Value = 1500

# Row 47
$ws.Cells.Item(47, 4).Value = 44692
$ws.Cells.Item(47, 10).Value = 100
$ws.Cells.Item(47, 11).Value = 25000
$ws.Cells.Item(47, 12).Value = 26000
$ws.Cells.Item(47, 13).Value = 25500
$ws.Cells.Item(47, 16).Value = 1020

# Row 48
$ws.Cells.Item(48, 4).Value = 44629
$ws.Cells.Item(48, 11).Value = 30000
$ws.Cells.Item(48, 12).Value = 30000
$ws.Cells.Item(48, 13).Value = 30000
$ws.Cells.Item(48, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(48, 16).Value = 1200

# Row 49
$ws.Cells.Item(49, 4).Value = 44489
$ws.Cells.Item(49, 11).Value = 40000
$ws.Cells.Item(49, 12).Value = 42000
$ws.Cells.Item(49, 13).Value = 41000
$ws.Cells.Item(49, 16).Value = 1640

# Row 50
$ws.Cells.Item(50, 4).Value = 44160
$ws.Cells.Item(50, 10).Value = 100
$ws.Cells.Item(50, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(50, 15).Value = 'Región de O''Higgins'

# Row 51
$ws.Cells.Item(51, 4).Value = 44433
$ws.Cells.Item(51, 11).Value = 25000
$ws.Cells.Item(51, 12).Value = 26000
$ws.Cells.Item(51, 13).Value = 25500
$ws.Cells.Item(51, 16).Value = 1020

# Row 52
$ws.Cells.Item(52, 4).Value = 44825
$ws.Cells.Item(52, 10).Value = 100
$ws.Cells.Item(52, 11).Value = 30000
$ws.Cells.Item(52, 12).Value = 32000
$ws.Cells.Item(52, 13).Value = 31000
$ws.Cells.Item(52, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(52, 15).Value = 'Perú'
$ws.Cells.Item(52, 16).Value = 1240

# Row 53
$ws.Cells.Item(53, 4).Value = 44294
$ws.Cells.Item(53, 10).Value = 100
$ws.Cells.Item(53, 11).Value = 24000
$ws.Cells.Item(53, 12).Value = 25000
$ws.Cells.Item(53, 13).Value = 24500
$ws.Cells.Item(53, 16).Value = 980

# Row 54
$ws.Cells.Item(54, 4).Value = 44244
$ws.Cells.Item(54, 11).Value = 16000
$ws.Cells.Item(54, 12).Value = 18000
$ws.Cells.Item(54, 13).Value = 17000
$ws.Cells.Item(54, 16).Value = 680

# Row 55
$ws.Cells.Item(55, 4).Value = 44384
$ws.Cells.Item(55, 8).Value = 'Sin especificar'
$ws.Cells.Item(55, 10).Value = 100
$ws.Cells.Item(55, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(55, 15).Value = 'Perú'

# Row 56
$ws.Cells.Item(56, 4).Value = 44881
$ws.Cells.Item(56, 11).Value = 45000
$ws.Cells.Item(56, 12).Value = 46000
$ws.Cells.Item(56, 13).Value = 45500
$ws.Cells.Item(56, 16).Value = 1820

# Row 58
$ws.Cells.Item(58, 4).Value = 44923
$ws.Cells.Item(58, 8).Value = 'Magnum'
$ws.Cells.Item(58, 11).Value = 27000
$ws.Cells.Item(58, 12).Value = 28000
$ws.Cells.Item(58, 13).Value = 27500
$ws.Cells.Item(58, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(58, 16).Value = 1100

# Row 59
$ws.Cells.Item(59, 4).Value = 44399
$ws.Cells.Item(59, 11).Value = 20000
$ws.Cells.Item(59, 12).Value = 22000
$ws.Cells.Item(59, 13).Value = 21000
$ws.Cells.Item(59, 16).Value = 840

# Row 60
$ws.Cells.Item(60, 4).Value = 44540
$ws.Cells.Item(60, 10).Value = 170
$ws.Cells.Item(60, 11).Value = 21000
$ws.Cells.Item(60, 12).Value = 22000
$ws.Cells.Item(60, 13).Value = 21529
$ws.Cells.Item(60, 16).Value = 861

# Row 61
$ws.Cells.Item(61, 4).Value = 44904
$ws.Cells.Item(61, 8).Value = 'Sin especificar'
$ws.Cells.Item(61, 10).Value = 350
$ws.Cells.Item(61, 11).Value = 42000
$ws.Cells.Item(61, 12).Value = 45000
$ws.Cells.Item(61, 13).Value = 43714
$ws.Cells.Item(61, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(61, 16).Value = 1749

# Row 62
$ws.Cells.Item(62, 4).Value = 44658
$ws.Cells.Item(62, 10).Value = 80
$ws.Cells.Item(62, 11).Value = 25000
$ws.Cells.Item(62, 13).Value = 25000
$ws.Cells.Item(62, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(62, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(62, 16).Value = 1000

# Row 63
$ws.Cells.Item(63, 4).Value = 44342
$ws.Cells.Item(63, 11).Value = 28000
$ws.Cells.Item(63, 12).Value = 30000
$ws.Cells.Item(63, 13).Value = 29000
$ws.Cells.Item(63, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(63, 16).Value = 1160

# Row 64
$ws.Cells.Item(64, 4).Value = 44855
$ws.Cells.Item(64, 11).Value = 30000
$ws.Cells.Item(64, 12).Value = 32000
$ws.Cells.Item(64, 13).Value = 31000
$ws.Cells.Item(64, 15).Value = 'Perú'
$ws.Cells.Item(64, 16).Value = 1240

# Row 65
$ws.Cells.Item(65, 4).Value = 44510
$ws.Cells.Item(65, 11).Value = 35000
$ws.Cells.Item(65, 12).Value = 36000
$ws.Cells.Item(65, 13).Value = 35500
$ws.Cells.Item(65, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(65, 15).Value = 'Perú'
$ws.Cells.Item(65, 16).Value = 1420

# Row 66
$ws.Cells.Item(66, 4).Value = 44636
$ws.Cells.Item(66, 10).Value = 180
$ws.Cells.Item(66, 11).Value = 22000
$ws.Cells.Item(66, 12).Value = 23000
$ws.Cells.Item(66, 13).Value = 22444
$ws.Cells.Item(66, 16).Value = 898

# Row 67
$ws.Cells.Item(67, 4).Value = 44272
$ws.Cells.Item(67, 11).Value = 22000
$ws.Cells.Item(67, 12).Value = 24000
$ws.Cells.Item(67, 13).Value = 23000
$ws.Cells.Item(67, 16).Value = 920


# --- New row 68 (new market price record, same block shape as existing rows) ---
$ws.Cells.Item(68, 1).Value = 11
$ws.Cells.Item(68, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(68, 3).Value = 'Bíobío'
$ws.Cells.Item(68, 4).Value = 44552
$ws.Cells.Item(68, 5).Value = 8
$ws.Cells.Item(68, 6).Value = 100112031
$ws.Cells.Item(68, 7).Value = 'Poroto verde'
$ws.Cells.Item(68, 8).Value = 'Magnum'
$ws.Cells.Item(68, 9).Value = 'Primera'
$ws.Cells.Item(68, 10).Value = 100
$ws.Cells.Item(68, 11).Value = 34000
$ws.Cells.Item(68, 12).Value = 36000
$ws.Cells.Item(68, 13).Value = 35000
$ws.Cells.Item(68, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(68, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(68, 16).Value = 1400
$ws.Cells.Item(68, 17).Value = 25
$ws.Cells.Item(68, 18).Value = 'Hortaliza'

# Column D (Fecha) carries the workbook's date-serial number format, same as the rest of the column
$ws.Cells.Item(68, 4).NumberFormat = $ws.Cells.Item(67, 4).NumberFormat
